$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7201
$ws.Range("C3").Value = 173693
$ws.Range("C5").Value = 9814
$ws.Range("C6").Value = 335
$ws.Range("C7").Value = 5.65
